$wb = $excel.ActiveWorkbook

# Rename sheet3: disc_list -> disk_list
$wsDisk = $wb.Worksheets.Item(3)
$wsDisk.Name = "disk_list"

# ===== Sheet 1: ip_address_list =====
$ws1 = $wb.Worksheets.Item("ip_address_list")
$ws1.Range("E1").Value = 0
$ws1.Range("B2").Value = "192.168.000.j"
$ws1.Range("D2").Value = "kkgg"
$ws1.Range("A3").Value = "bewolktEN"
$ws1.Range("B3").Value = "192.168.000.000"
$ws1.Range("C3").Value = "255.255.255.0"
$ws1.Range("D3").Value = "du hast einen problem"
$ws1.Range("E3").Value = 1
$ws1.Range("A4").Value = "einkaufenfh"
$ws1.Range("B4").Value = "192.168.000.000"
$ws1.Range("C4").Value = "255.255.255.0"
$ws1.Range("D4").Value = "gggg"
$ws1.Range("E4").Value = 0

# ===== Sheet 2: ip_adress_fav_list =====
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")
$ws2.Range("A1").Value = "bewolktEN"
$ws2.Range("D1").Value = "du hast einen problem"

# ===== Sheet 4: Settings =====
$ws4 = $wb.Worksheets.Item("Settings")
$ws4.Range("B2").Value = "Ethernet,Ethernet 1,Ethernet 2,Ethernet 3,Ethernet 4,Ethernet 5,Wi-Fi,"
$ws4.Range("A3").Value = "spousteci okno: na oblibenych = 1"
$ws4.Range("A4").Value = "spousteci okno: zobrazeni disku = 1, normal = 0"
$ws4.Range("B4").Value = 1
$ws4.Range("A5").Value = "zakladni velikost okna normal = 0, max = 1,min = 2"
$ws4.Range("B5").Value = 2

# ===== Selections (set last so final active sheet/tab matches the target) =====
# Sheet1 selection
$ws1.Range("G13").Select()
# Sheet2 selection
$ws2.Range("E8").Select()
# Sheet3 (disk_list) selection
$wsDisk.Range("G33").Select()
# Sheet4 (Settings) selection - Settings ends up the active sheet/tab
$ws4.Select()
$ws4.Range("F18").Select()
